$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing ISBN for row 2
$ws.Range("A2").Value = "111-11-11111-11-1"

# Move the active selection to A3 (as left after the edit)
$ws.Range("A3").Select()
